$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new columns at Q (17) so the old Remarks/file_name columns
# --- shift right to become S (Remarks) and T (file_name). ---
$ws.Columns.Item(17).Insert()
$ws.Columns.Item(17).Insert()

# --- Header row: label the two newly-inserted columns ---
$ws.Range("Q1").Value = "Phase_Code"
$ws.Range("R1").Value = "Cost_Type"

# --- Column widths for the changed / new columns ---
$ws.Columns.Item(5).ColumnWidth = 14
$ws.Columns.Item(7).ColumnWidth = 16
$ws.Columns.Item(17).ColumnWidth = 12
$ws.Columns.Item(18).ColumnWidth = 11
$ws.Columns.Item(19).ColumnWidth = 9
$ws.Columns.Item(20).ColumnWidth = 32

# --- Row 2: updated vendor / invoice values ---
$ws.Range("D2").Value = "CAPAIR"
$ws.Range("E2").Value = "Captive Aire"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "2457243"
$ws.Range("I2:L2").NumberFormat = "@"
$ws.Range("I2").Value = "4852929"
$ws.Range("J2").Value = "441685"
$ws.Range("K2").Value = "102125"
$ws.Range("L2").Value = "4309119"
$ws.Range("M2").Value = 1412
$ws.Range("P2").NumberFormat = "@"
$ws.Range("P2").Value = "5030"

# --- Row 2: newly-inserted Phase_Code / Cost_Type values ---
$ws.Range("Q2").Value = 320
$ws.Range("R2").Value = "M"

# --- Row 2: Remarks (now S2) no longer carries a value ---
$ws.Range("S2").Value = ""

# --- Row 2: file_name (now T2) ---
$ws.Range("T2").Value = "captive aire_1754662125633.pdf"
